$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Rewrite the Abstract paragraph + the blank paragraph that follows it.
#    Old layout:
#      Para A (style=Abstract): "Abstract" + "—" + Spanish placeholder text
#      Para B (blank, rPr lang=es-CO)
#    New layout:
#      Para A (no pStyle / default): "Abstract" + "—" + " " + English abstract text (2 runs)
#      Para B (style=Abstract, rFonts Arial, no lang): blank
#      Para C: brand new, fully empty paragraph
# ------------------------------------------------------------------
$abstractPara = $d.Paragraphs(7)
$blankPara = $d.Paragraphs(8)
$spanRange = $d.Range($abstractPara.Range.Start, $blankPara.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:color w:val="0070C0"/></w:rPr><w:t>Abstract</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="0070C0"/></w:rPr><w:t>&#8212;</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>This project aims to explore the strengths and weaknesses of centrality indices when applied to affiliation networks. The case study involves examining the affiliation network of corporate executive officers in Colombia and their membership in various services such as clubs, health service providers, and recreation service providers. The dataset used for this study consists of membership information of the corporate executives in social organizations, generating a bipartite network where left nodes represent persons and right nodes represent social organizations.</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>The implemented network science approach involves generating a bipartite affiliation network and analyzing the centrality of the network. Additionally, a visual representation of the network will be created from the dataset. This study sheds light on the use of centrality measures in affiliation networks and their potential strengths and limitations.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Abstract"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>
</w:p>
<w:p/>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$spanRange.InsertXML($xml)

# ------------------------------------------------------------------
# 2) Split the "Goal And Scope" paragraph run into two runs at
#    "...of centrality " | "indices when applied..."
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "This project discusses strengths and weaknesses of centrality indices when applied to affiliation networks.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This project discusses strengths and weaknesses of centrality indices when applied to affiliation networks.",
    2)

# ------------------------------------------------------------------
# 3) Merge the two runs in the "Affiliation Networks" paragraph that
#    were split around "...subset of " | "events and each event...".
#    Removing the mid-sentence split simply requires deleting the
#    run boundary; since both runs carry identical formatting, a
#    plain text find/replace collapses them into a single run.
# ------------------------------------------------------------------
$found2 = $d.Content.Find.Execute(
    "relates each actor to a subset of events and each event to a subset of actors.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "relates each actor to a subset of events and each event to a subset of actors.",
    2)
